# Commit message: "Change Excel Field View to Cache, And set default value to FALSE"
#
# The "Property" sheet has a header row (row 1) describing field metadata
# columns. Column F's header was "View" and its per-row default values
# (rows 2-25) were all TRUE. Rename the header to "Cache" and flip every
# default value in that column to FALSE.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property")

# Rename column F's header from "View" to "Cache".
$ws.Range("F1").Value = "Cache"

# Flip the default value for every data row in column F from TRUE to FALSE.
$ws.Range("F2:F25").Value = $false

# Reflect the edited column as the active selection, matching the
# author's on-screen state after making the change.
$ws.Range("F2:F25").Select() | Out-Null
